# Generate Report for Archive
# Update the localization status from "Ready for handoff" to "In Translation"
# across the Overview sheet and each locale sheet, and resize the now-narrower
# status columns to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status text ---
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

# --- Locale detail sheets: column C holds the Status text ---
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# --- Narrow the status columns now that the text is shorter ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede.Columns.Item(3).ColumnWidth = 12.5
